
$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

# --- Sanity checks: verify expected occurrence counts before mutating ---
$old1 = '<w:r w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t>Télécom SudParis – Étudiante ingénieure</w:t></w:r>'
$new1 = '<w:r w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Télécom </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t>SudParis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> – Étudiante I</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t>ngénieure</w:t></w:r>'
$old2 = '<w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00B34C7C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t>- Licence en Ingénierie</w:t></w:r><w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r>'
$new2id0 = '<w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B34C7C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Licence en Ingénierie</w:t></w:r><w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r>'
$new2id1 = '<w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B34C7C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="193C61"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Licence en Ingénierie</w:t></w:r><w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r w:rsidR="008D798C" w:rsidRPr="00C6036E"><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r>'
$old3 = '<w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Réaliser un jeu </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>Sudoku avec interfaces graphiques en utilisant le langage JAVA et l’outil SVN</w:t></w:r>'
$new3 = '<w:r><w:rPr><w:rFonts w:eastAsia="Microsoft YaHei" w:cstheme="minorHAnsi"/><w:color w:val="5B5B5B"/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>Réaliser un jeu Sudoku avec interfaces graphiques en utilisant le langage JAVA et l’outil SVN</w:t></w:r>'
$old4 = '"><v:shape id="_x0000_s1035"'
$new4 = '"><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s1035"'

function CountOccurrences($haystack, $needle) {
    $count = 0
    $idx = 0
    while ($true) {
        $pos = $haystack.IndexOf($needle, $idx)
        if ($pos -lt 0) { break }
        $count += 1
        $idx = $pos + $needle.Length
    }
    return $count
}

Write-Output "old1 count:"
Write-Output (CountOccurrences $xml $old1)
Write-Output "old2 count:"
Write-Output (CountOccurrences $xml $old2)
Write-Output "old3 count:"
Write-Output (CountOccurrences $xml $old3)
Write-Output "old4 count:"
Write-Output (CountOccurrences $xml $old4)

# --- 1) Télécom SudParis run-split (applies identically to both occurrences) ---
$xml = $xml.Replace($old1, $new1)

# --- 2) "- Licence en Ingénierie" block: two occurrences, need distinct bookmark ids (0 then 1) ---
$pos1 = $xml.IndexOf($old2)
if ($pos1 -ge 0) {
    $before = $xml.Substring(0, $pos1)
    $after = $xml.Substring($pos1 + $old2.Length)
    $pos2 = $after.IndexOf($old2)
    if ($pos2 -ge 0) {
        $middle = $after.Substring(0, $pos2)
        $tail = $after.Substring($pos2 + $old2.Length)
        $xml = $before + $new2id0 + $middle + $new2id1 + $tail
    } else {
        $xml = $before + $new2id0 + $after
    }
}

# --- 3) Merge "Réaliser un jeu" + "Sudoku avec..." runs, drop the old bookmark (id=0 occurrence only) ---
$xml = $xml.Replace($old3, $new3)

# --- 4) Insert v:shapetype boilerplate before the VML textbox shape ---
$xml = $xml.Replace($old4, $new4)

$d.Content.WordOpenXML = $xml
Write-Output "DONE"
